# fix: standardize "Notice u/s 94 BNSS, 2023" line above "To," across the
# template - replace the messy, copy/pasted paragraph (inherited "Body A"
# style, center alignment, explicit Bookman-Old-Style/size/color run
# formatting) with a plain bold+underline run in a justified paragraph.

$d = $word.ActiveDocument
$targetText = "Notice u/s 94 BNSS, 2023"

# Minimal OOXML package payload used to splice in a clean paragraph via
# Range.InsertXML - this lets us land *exactly* the desired <w:pPr>/<w:rPr>
# (no leftover pStyle/shd/spacing/rFonts/color/sz/etc.) instead of the
# residue that setting individual Font/ParagraphFormat properties leaves
# behind.
$cleanParaXmlTemplate = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:pPr><w:jc w:val="both"/></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>{0}</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$cleanParaXml = $cleanParaXmlTemplate -replace '\{0\}', $targetText

# Walk every occurrence of the line (the same template text can appear more
# than once across the document) and standardize each one.
$searchRng = $d.Content.Duplicate
$searchRng.Start = 0
$searchRng.End = $d.Content.End

while ($true) {
    [void]$searchRng.Find.Execute($targetText, $true, $true, $false, $false, $false,
                                   $true, 0, $false, "", 0)
    if (-not $searchRng.Find.Found) { break }

    $para = $searchRng.Paragraphs.First
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End

    # Insert the freshly-formatted paragraph right after the existing one
    # (collapsed range at the paragraph end lands it as its own new <w:p>,
    # sibling to the old paragraph, carrying only the pPr/rPr we specify).
    $insertionPoint = $d.Range($pEnd, $pEnd)
    [void]$insertionPoint.InsertXML($cleanParaXml)

    # Remove the original, messily-formatted paragraph (including its
    # paragraph mark) now that its replacement sits right after it.
    $oldPara = $d.Range($pStart, $pEnd + 1)
    [void]$oldPara.Delete()

    # Resume searching after the text we just normalized.
    $searchRng.Start = $pStart + $targetText.Length
    $searchRng.End = $d.Content.End
}
